$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(29, "test again", "04-11-2023"),
    @(30, "rewr", "04-11-2023"),
    @(31, "sawewe", "04-11-2023"),
    @(32, "erwerewrw", "04-11-2023"),
    @(33, "ewrwer", "04-11-2023"),
    @(34, "wewqe", "04-11-2023"),
    @(35, "dfsfsdfsfsd", "04-11-2023")
)

$startRow = 31
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]

    # Column C holds date-like text (e.g. "04-11-2023") that must stay text,
    # not get auto-converted into a date serial. Force text via "@" format,
    # then restore the default "Normal" style so no extra per-cell style
    # index lingers on the written cell.
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 3).Style = "Normal"
}
